$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.345.13'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '3.743.36'
$ws.Range('E3').Value = '  -1.05%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'591.99"
$ws.Range('E5').Value = '  -0.80%  '
$ws.Range('D6').Value = "'165.36"
$ws.Range('E6').Value = '  -1.79%  '
$ws.Range('D7').Value = '3.741.57'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -1.43%  '
$ws.Range('E10').Value = '  -3.14%  '
$ws.Range('E11').Value = '  -1.76%  '
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').Value = "'0.0000258"
$ws.Range('E13').Value = '  -7.33%  '
$ws.Range('D14').Value = "'35.83"
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').Value = '4.369.95'
$ws.Range('E15').Value = '  -1.05%  '
$ws.Range('D16').Value = '3.760.27'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('D17').Value = '68.275.58'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('E18').Value = '  -4.27%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = "'0.111"
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = "'6.94"
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('E21').Value = '  +1.15%  '
$ws.Range('D22').Value = "'462.42"
$ws.Range('E22').Value = '  -1.04%  '
$ws.Range('E23').Value = '  -3.54%  '
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('E25').Value = '  -2.96%  '
$ws.Range('E26').Value = '  -3.95%  '
$ws.Range('D27').Value = "'11.82"
$ws.Range('E27').Value = '  -2.88%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = "'9.98"
$ws.Range('E29').Value = '  -4.04%  '
$ws.Range('D30').Value = '3.889.76'
$ws.Range('E30').Value = '  -1.01%  '
$ws.Range('D31').Value = "'2.75"
$ws.Range('E31').Value = '  -5.68%  '
$ws.Range('D32').Value = "'7.27"
$ws.Range('E32').Value = '  -4.12%  '
$ws.Range('D33').Value = "'29.75"
$ws.Range('E33').Value = '  -2.30%  '
$ws.Range('D34').Value = "'2.14"
$ws.Range('E34').Value = '  -3.56%  '
$ws.Range('D35').Value = "'0.999"
$ws.Range('D36').Value = "'9.09"
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('D37').Value = '3.698.60'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').Value = "'0.0998"
$ws.Range('E38').Value = '  -3.84%  '
$ws.Range('E39').Value = '  -10.75%  '
$ws.Range('D40').Value = "'1.00"
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('E42').Value = '  -1.90%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E45').Value = '  -4.18%  '
$ws.Range('D46').Value = "'43.10"
$ws.Range('E46').Value = '  +7.75%  '
$ws.Range('E47').Value = '  +2.14%  '
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('D50').Value = "'144.72"
$ws.Range('E50').Value = '  +0.59%  '
$ws.Range('D51').Value = "'386.09"
$ws.Range('E51').Value = '  -4.90%  '
